$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21 (Anna / 88. Merge Sorted Array): remove the incomplete date in E21 ---
$ws.Range("E21").Clear()

# --- Row 23 (Anna / 176. Second Highest Salary): remove the incomplete date in E23 ---
$ws.Range("E23").Clear()

# --- Row 25 (Anna / 234. Palindrome Linked List) ---
$ws.Range("D25").WrapText = $true
$ws.Range("D25").Value = "234. Palindrome Linked List"
$ws.Rows.Item(25).RowHeight = 14

# --- Row 26 (Stephan / 234. Palindrome Linked List) ---
$ws.Range("D26").WrapText = $true
$ws.Range("D26").Value = "234. Palindrome Linked List"
$ws.Range("E26").Value = "2020/12/15"
$ws.Range("F26").Value = "Lindked List"
$ws.Range("G26").Value = "Completed"
$ws.Rows.Item(26).RowHeight = 14

# --- move the active selection like the author's session ended up ---
$ws.Range("F33").Select()
